$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'81.409.83"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").Value = "'3.148.74"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'209.25"
$ws.Range("E5").Value = "  +2.15%  "
$ws.Range("D6").Value = "'620.66"
$ws.Range("E6").Value = "  -2.19%  "
$ws.Range("D7").Value = "'0.280"
$ws.Range("E7").Value = "  +21.98%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.581"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "'3.147.22"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").Value = "'0.581"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "'0.0000249"
$ws.Range("E12").Value = "  +10.41%  "
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").Value = "'5.27"
$ws.Range("E14").Value = "  -3.97%  "
$ws.Range("D15").Value = "'3.726.35"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").Value = "'31.32"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "'81.034.64"
$ws.Range("E17").Value = "  +2.64%  "
$ws.Range("D18").Value = "'3.141.51"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").Value = "'3.16"
$ws.Range("E19").Value = "  +4.17%  "
$ws.Range("D20").Value = "'13.85"
$ws.Range("E20").Value = "  -4.31%  "
$ws.Range("D21").Value = "'430.12"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").Value = "'8.90"
$ws.Range("E22").Value = "  -3.69%  "
$ws.Range("D23").Value = "'5.09"
$ws.Range("E23").Value = "  +1.86%  "
$ws.Range("D24").Value = "'7.23"
$ws.Range("E24").Value = "  +5.28%  "
$ws.Range("E25").Value = "  +9.12%  "
$ws.Range("D26").Value = "'3.312.03"
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("D27").Value = "'76.22"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").Value = "'10.75"
$ws.Range("E28").Value = "  -3.72%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'0.0000120"
$ws.Range("E30").Value = "  +3.19%  "
$ws.Range("D31").Value = "'580.40"
$ws.Range("E31").Value = "  +11.25%  "
$ws.Range("D32").Value = "'0.998"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").Value = "'8.91"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "'1.50"
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("D35").Value = "'0.153"
$ws.Range("E35").Value = "  +10.11%  "
$ws.Range("D36").Value = "'0.140"
$ws.Range("E36").Value = "  +13.85%  "
$ws.Range("D37").Value = "'1.98"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").Value = "'22.65"
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("D39").Value = "'0.998"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").Value = "'6.01"
$ws.Range("E40").Value = "  +10.32%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.405"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").Value = "'2.06"
$ws.Range("E42").Value = "  +14.87%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "'20.72"
$ws.Range("E43").Value = "  +3.58%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'3.00"
$ws.Range("E44").Value = "  +19.44%  "
$ws.Range("D45").Value = "'158.44"
$ws.Range("E45").Value = "  -3.72%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'186.25"
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("D48").Value = "'45.22"
$ws.Range("E48").Value = "  +5.92%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'0.770"
$ws.Range("E50").Value = "  -5.40%  "
$ws.Range("D51").Value = "'25.90"
$ws.Range("E51").Value = "  +0.97%  "
